$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the word list (shared strings) for the tail of the existing
#     sentence block (C170 = block 8), rows 186-189, columns E (word) / F (count) ---
$ws.Range("E186").Value = "[b'make']"
$ws.Range("F186").Value = 1

$ws.Range("E187").Value = "[b'any']"
$ws.Range("F187").Value = 1

$ws.Range("E188").Value = "[b'disciplinary']"
$ws.Range("F188").Value = 1

$ws.Range("E189").Value = "[b'decisions']"
$ws.Range("F189").Value = 1

# --- Append a brand-new (still-empty) sentence block: C191 = 9 (block index),
#     D191:D210 = 0..19 (per-sentence word index), no words filled in yet ---
$ws.Range("C191").Value = 9
$ws.Range("D191").Value = 0

$row = 192
for ($i = 1; $i -le 19; $i++) {
    $ws.Range("D$row").Value = $i
    $row = $row + 1
}

# --- Update the saved view state to match where Excel was scrolled/selected ---
$null = $ws.Range("E202").Select()
$excel.ActiveWindow.ScrollRow = 180
$excel.ActiveWindow.ScrollColumn = 1
